$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The single "2030" scenario row (row 6) is being expanded into a full
# set of milestone years (2025, 2030, 2035, 2040, 2045, 2050), each as
# its own row, by copying the original data row (with its formatting)
# down and then updating the Year value (column F) in each copy.

# First, change the existing row 6 value from 2030 to 2025 (it becomes
# the first entry of the new year series).
$ws.Range("F6").Value2 = 2025

# Duplicate row 6 (B6:H6, including formatting/styles) into rows 7-11.
$ws.Range("B6:H6").Copy($ws.Range("B7:H7"))
$ws.Range("B6:H6").Copy($ws.Range("B8:H8"))
$ws.Range("B6:H6").Copy($ws.Range("B9:H9"))
$ws.Range("B6:H6").Copy($ws.Range("B10:H10"))
$ws.Range("B6:H6").Copy($ws.Range("B11:H11"))

# Set the Year value for each newly added row.
$ws.Range("F7").Value2 = 2030
$ws.Range("F8").Value2 = 2035
$ws.Range("F9").Value2 = 2040
$ws.Range("F10").Value2 = 2045
$ws.Range("F11").Value2 = 2050

# Update the active selection to reflect where the user ended up
# after editing (below the newly added table).
$ws.Range("E16").Select()
